# Normalize the "Recorded By" list in column G: the exact token "System"
# (capital S) is moved to the end of the comma-separated list, preserving
# the relative order of the remaining entries. When the list has no
# "System" token at all, the two entries simply swap places (list reversed).
# Lists with a single entry are left untouched.

function Transform-RecordedBy {
    param([string]$s)

    $parts = $s.Split(",") | ForEach-Object { $_.Trim() }

    if ($parts.Count -le 1) {
        return $s
    }

    $rest = @()
    $sys = @()
    foreach ($p in $parts) {
        if ($p.Equals("System")) {
            $sys += $p
        } else {
            $rest += $p
        }
    }

    if ($sys.Count -gt 0) {
        $newParts = $rest + $sys
    } else {
        $newParts = @()
        for ($i = $parts.Count - 1; $i -ge 0; $i--) {
            $newParts += $parts[$i]
        }
    }

    return ($newParts -join ", ")
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $current = $cell.Value()
    if ($current -ne $null -and $current -ne "") {
        $updated = Transform-RecordedBy $current
        if ($updated -ne $current) {
            $cell.Value = $updated
        }
    }
}
